$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Tgfb3"
$ws.Cells.Item(2, 3).Value = "Tgfbr1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.6423476666666667
$ws.Cells.Item(2, 8).Value = 1.927043
$ws.Cells.Item(2, 9).Value = 0.01173234890143342
$ws.Cells.Item(2, 10).Value = 0.01173234890143342
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 68.65869266666667
$ws.Cells.Item(2, 14).Value = 205.976078
$ws.Cells.Item(2, 15).Value = 0.6475952735309433
$ws.Cells.Item(2, 16).Value = 0.6475952735309431
$ws.Cells.Item(2, 17).Value = 44.10275103081711
$ws.Cells.Item(2, 18).Value = 396.924759277354
$ws.Cells.Item(2, 19).Value = 0.007597813695984237
$ws.Cells.Item(2, 20).Value = 0.007597813695984237

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Tgfb3"
$ws.Cells.Item(3, 3).Value = "Tgfbr1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.6423476666666667
$ws.Cells.Item(3, 8).Value = 1.927043
$ws.Cells.Item(3, 9).Value = 0.01173234890143342
$ws.Cells.Item(3, 10).Value = 0.01173234890143342
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 28.23943666666667
$ws.Cells.Item(3, 14).Value = 84.71831
$ws.Cells.Item(3, 15).Value = 0.2663570336431459
$ws.Cells.Item(3, 16).Value = 0.2663570336431459
$ws.Cells.Item(3, 17).Value = 18.13953625081444
$ws.Cells.Item(3, 18).Value = 163.25582625733
$ws.Cells.Item(3, 19).Value = 0.003124993651052227
$ws.Cells.Item(3, 20).Value = 0.003124993651052227

$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Tgfb3"
$ws.Cells.Item(4, 3).Value = "Tgfbr1"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.6423476666666667
$ws.Cells.Item(4, 8).Value = 1.927043
$ws.Cells.Item(4, 9).Value = 0.01173234890143342
$ws.Cells.Item(4, 10).Value = 0.01173234890143342
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 9.122861666666667
$ws.Cells.Item(4, 14).Value = 27.368585
$ws.Cells.Item(4, 15).Value = 0.08604769282591093
$ws.Cells.Item(4, 16).Value = 0.08604769282591092
$ws.Cells.Item(4, 17).Value = 5.860048904906111
$ws.Cells.Item(4, 18).Value = 52.740440144155
$ws.Cells.Item(4, 19).Value = 0.001009541554396956
$ws.Cells.Item(4, 20).Value = 0.001009541554396956

$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Tgfb3"
$ws.Cells.Item(5, 3).Value = "Tgfbr1"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 22.27635266666667
$ws.Cells.Item(5, 8).Value = 66.829058
$ws.Cells.Item(5, 9).Value = 0.4068730304461968
$ws.Cells.Item(5, 10).Value = 0.4068730304461968
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 68.65869266666667
$ws.Cells.Item(5, 14).Value = 205.976078
$ws.Cells.Item(5, 15).Value = 0.6475952735309433
$ws.Cells.Item(5, 16).Value = 0.6475952735309431
$ws.Cells.Item(5, 17).Value = 1529.465251474947
$ws.Cells.Item(5, 18).Value = 13765.18726327453
$ws.Cells.Item(5, 19).Value = 0.2634890514441686
$ws.Cells.Item(5, 20).Value = 0.2634890514441686

$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Tgfb3"
$ws.Cells.Item(6, 3).Value = "Tgfbr1"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 22.27635266666667
$ws.Cells.Item(6, 8).Value = 66.829058
$ws.Cells.Item(6, 9).Value = 0.4068730304461968
$ws.Cells.Item(6, 10).Value = 0.4068730304461968
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 28.23943666666667
$ws.Cells.Item(6, 14).Value = 84.71831
$ws.Cells.Item(6, 15).Value = 0.2663570336431459
$ws.Cells.Item(6, 16).Value = 0.2663570336431459
$ws.Cells.Item(6, 17).Value = 629.0716502946644
$ws.Cells.Item(6, 18).Value = 5661.644852651981
$ws.Cells.Item(6, 19).Value = 0.1083734934590464
$ws.Cells.Item(6, 20).Value = 0.1083734934590464

$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Tgfb3"
$ws.Cells.Item(7, 3).Value = "Tgfbr1"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 22.27635266666667
$ws.Cells.Item(7, 8).Value = 66.829058
$ws.Cells.Item(7, 9).Value = 0.4068730304461968
$ws.Cells.Item(7, 10).Value = 0.4068730304461968
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 9.122861666666667
$ws.Cells.Item(7, 14).Value = 27.368585
$ws.Cells.Item(7, 15).Value = 0.08604769282591093
$ws.Cells.Item(7, 16).Value = 0.08604769282591092
$ws.Cells.Item(7, 17).Value = 203.2240838158811
$ws.Cells.Item(7, 18).Value = 1829.01675434293
$ws.Cells.Item(7, 19).Value = 0.03501048554298185
$ws.Cells.Item(7, 20).Value = 0.03501048554298184

$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Tgfb3"
$ws.Cells.Item(8, 3).Value = "Tgfbr1"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 31.831433
$ws.Cells.Item(8, 8).Value = 95.494299
$ws.Cells.Item(8, 9).Value = 0.5813946206523698
$ws.Cells.Item(8, 10).Value = 0.5813946206523697
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 68.65869266666667
$ws.Cells.Item(8, 14).Value = 205.976078
$ws.Cells.Item(8, 15).Value = 0.6475952735309433
$ws.Cells.Item(8, 16).Value = 0.6475952735309431
$ws.Cells.Item(8, 17).Value = 2185.504575486591
$ws.Cells.Item(8, 18).Value = 19669.54117937932
$ws.Cells.Item(8, 19).Value = 0.3765084083907904
$ws.Cells.Item(8, 20).Value = 0.3765084083907903

$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Tgfb3"
$ws.Cells.Item(9, 3).Value = "Tgfbr1"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 31.831433
$ws.Cells.Item(9, 8).Value = 95.494299
$ws.Cells.Item(9, 9).Value = 0.5813946206523698
$ws.Cells.Item(9, 10).Value = 0.5813946206523697
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 28.23943666666667
$ws.Cells.Item(9, 14).Value = 84.71831
$ws.Cells.Item(9, 15).Value = 0.2663570336431459
$ws.Cells.Item(9, 16).Value = 0.2663570336431459
$ws.Cells.Item(9, 17).Value = 898.9017362127433
$ws.Cells.Item(9, 18).Value = 8090.11562591469
$ws.Cells.Item(9, 19).Value = 0.1548585465330473
$ws.Cells.Item(9, 20).Value = 0.1548585465330473

$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Tgfb3"
$ws.Cells.Item(10, 3).Value = "Tgfbr1"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 31.831433
$ws.Cells.Item(10, 8).Value = 95.494299
$ws.Cells.Item(10, 9).Value = 0.5813946206523698
$ws.Cells.Item(10, 10).Value = 0.5813946206523697
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 9.122861666666667
$ws.Cells.Item(10, 14).Value = 27.368585
$ws.Cells.Item(10, 15).Value = 0.08604769282591093
$ws.Cells.Item(10, 16).Value = 0.08604769282591092
$ws.Cells.Item(10, 17).Value = 290.3937599107683
$ws.Cells.Item(10, 18).Value = 2613.543839196915
$ws.Cells.Item(10, 19).Value = 0.05002766572853213
$ws.Cells.Item(10, 20).Value = 0.05002766572853211
Write-Host "Edit applied"
